$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.702.71"
$ws.Range("E2").Value = "  +5.71%  "
$ws.Range("D3").Value = "3.355.20"
$ws.Range("E3").Value = "  +10.53%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'256.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.02%  "
$ws.Range("D6").Value = "'622.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("D7").Value = "'1.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.38%  "
$ws.Range("D8").Value = "'0.386"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.38%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.353.38"
$ws.Range("E10").Value = "  +10.61%  "
$ws.Range("D11").Value = "'0.803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "'0.199"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").Value = "98.374.86"
$ws.Range("E13").Value = "  +5.65%  "
$ws.Range("D14").Value = "'35.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.95%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.984.42"
$ws.Range("E15").Value = "  +10.54%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000245"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.58%  "
$ws.Range("D17").Value = "'5.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.46%  "
$ws.Range("D18").Value = "3.360.52"
$ws.Range("E18").Value = "  +10.04%  "
$ws.Range("D19").Value = "'3.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("D20").Value = "'15.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.95%  "
$ws.Range("D21").Value = "'483.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.47%  "
$ws.Range("D22").Value = "'5.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.30%  "
$ws.Range("D23").Value = "'0.0000206"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.92%  "
$ws.Range("D24").Value = "'9.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.55%  "
$ws.Range("D25").Value = "'5.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.93%  "
$ws.Range("D26").Value = "'88.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.52%  "
$ws.Range("D27").Value = "'11.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.45%  "
$ws.Range("D28").Value = "3.541.31"
$ws.Range("E28").Value = "  +10.78%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").Value = "'0.187"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.98%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.249"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.47%  "
$ws.Range("D32").Value = "'0.127"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.84%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.59%  "
$ws.Range("D34").Value = "'9.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("D35").Value = "'27.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.77%  "
$ws.Range("D36").Value = "'7.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "'517.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.80%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'1.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.93%  "
$ws.Range("D40").Value = "'24.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("D41").Value = "'0.446"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("E42").Value = "  +3.81%  "
$ws.Range("D43").Value = "'3.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("D44").Value = "'3.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.33%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'0.773"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.97%  "
$ws.Range("D47").Value = "'160.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'1.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.51%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'45.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.47%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.19%  "
$ws.Range("D51").Value = "'4.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.83%  "
